$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '66.799.05'
$ws.Range('D2').Style = "Normal"

$ws.Range('E2').Value = '  +2.07%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.093.32'
$ws.Range('D3').Style = "Normal"

$ws.Range('E3').Value = '  +5.27%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '580.28'
$ws.Range('D5').Style = "Normal"

$ws.Range('E5').Value = '  +2.06%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '168.62'
$ws.Range('D6').Style = "Normal"

$ws.Range('E6').Value = '  +6.18%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = "Normal"

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.087.72'
$ws.Range('D8').Style = "Normal"

$ws.Range('E8').Value = '  +5.19%  '

$ws.Range('E9').Value = '  +1.68%  '

$ws.Range('E10').Value = '  -2.07%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.156'
$ws.Range('D11').Style = "Normal"

$ws.Range('E11').Value = '  +3.77%  '

$ws.Range('E12').Value = '  +4.78%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000251'
$ws.Range('D13').Style = "Normal"

$ws.Range('E13').Value = '  +2.57%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '36.43'
$ws.Range('D14').Style = "Normal"

$ws.Range('E14').Value = '  +6.03%  '

$ws.Range('E15').Value = '  -0.58%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.605.08'
$ws.Range('D16').Style = "Normal"

$ws.Range('E16').Value = '  +5.26%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '66.748.88'
$ws.Range('D17').Style = "Normal"

$ws.Range('E17').Value = '  +2.03%  '

$ws.Range('E18').Value = '  +2.79%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.092.58'
$ws.Range('D19').Style = "Normal"

$ws.Range('E19').Value = '  +5.28%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '16.23'
$ws.Range('D20').Style = "Normal"

$ws.Range('E20').Value = '  +3.36%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '467.11'
$ws.Range('D21').Style = "Normal"

$ws.Range('E21').Value = '  +5.04%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.715'
$ws.Range('D22').Style = "Normal"

$ws.Range('E22').Value = '  +3.35%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.52'
$ws.Range('D23').Style = "Normal"

$ws.Range('E23').Value = '  +3.19%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '84.00'
$ws.Range('D24').Style = "Normal"

$ws.Range('E25').Value = '  +4.98%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '13.08'
$ws.Range('D26').Style = "Normal"

$ws.Range('E26').Value = '  +7.90%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.12'
$ws.Range('D27').Style = "Normal"

$ws.Range('E27').Value = '  +0.60%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.04'
$ws.Range('D29').Style = "Normal"

$ws.Range('E29').Value = '  +0.07%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.41'
$ws.Range('D30').Style = "Normal"

$ws.Range('E30').Value = '  +1.93%  '

$ws.Range('E32').Value = '  +1.17%  '

$ws.Range('E33').Value = '  +4.48%  '

$ws.Range('E34').Value = '  +3.49%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = "Normal"

$ws.Range('E35').Value = '  +0.00%  '

$ws.Range('E36').Value = '  +3.37%  '

$ws.Range('E37').Value = '  +2.74%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '47.34'
$ws.Range('D38').Style = "Normal"

$ws.Range('E38').Value = '  +5.43%  '

$ws.Range('E39').Value = '  +6.28%  '

$ws.Range('E40').Value = '  +6.13%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '50.30'
$ws.Range('D41').Style = "Normal"

$ws.Range('E41').Value = '  +1.24%  '

$ws.Range('E42').Value = '  +1.79%  '

$ws.Range('E43').Value = '  +2.44%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.83'
$ws.Range('D44').Style = "Normal"

$ws.Range('E44').Value = '  -0.31%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '382.63'
$ws.Range('D46').Style = "Normal"

$ws.Range('E46').Value = '  -0.17%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.789.61'
$ws.Range('D47').Style = "Normal"

$ws.Range('E47').Value = '  +3.36%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '134.88'
$ws.Range('D48').Style = "Normal"

$ws.Range('E48').Value = '  +1.08%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '24.97'
$ws.Range('D50').Style = "Normal"

$ws.Range('E50').Value = '  +6.90%  '

$ws.Range('E51').Value = '  +1.55%  '
